# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (spreadsheet column G, header "K") is recalculated from the
# updated data pipeline and its 64 data-row values (rows 2-65) are rewritten
# here with the freshly computed results.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @(0,2,3,2,0,2,0,1,1,1,3,1,1,2,1,2,0,4,7,1,0,1,2,1,2,3,1,4,3,1,5,1,1,2,9,1,2,2,3,3,3,0,2,3,3,2,2,3,2,3,1,1,2,1,0,3,5,1,3,3,1,3,2,1)

$startRow = 2
for ($i = 0; $i -lt $newK.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $newK[$i]
}
